# Update dSF (column F) values to reflect repulled data / recalculated mean
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F4").Value = -5
$ws.Range("F5").Value = -6
$ws.Range("F7").Value = 0
$ws.Range("F16").Value = -2
$ws.Range("F20").Value = 5
